$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.017
$ws.Range("E2").Value = 0.0258
$ws.Range("G2").Value = 0.005280456921170322
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2.444
$ws.Range("L2").Value = 0.13168812974837
$ws.Range("M2").Value = 1.77
$ws.Range("N2").Value = 0.004696205890156541
$ws.Range("O2").Value = 0.7242225859247136
$ws.Range("P2").Value = 1.77
$ws.Range("Q2").Value = 0.004696205890156541
$ws.Range("R2").Value = 0.7242225859247136
$ws.Range("U2").Value = 116.906
$ws.Range("V2").Value = 0.3101777659856726
$ws.Range("W2").Value = -0.05333333333333333
$ws.Range("X2").Value = 0.03673057250252455
$ws.Range("Y2").Value = -0.09006390583585788
$ws.Range("Z2").Value = -0.3739698148185466
$ws.Range("AB2").Value = 0.03669751335046647
$ws.Range("AC2").Value = -0.03669751335046647
$ws.Range("AD2").Value = 92.27
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 92.27
$ws.Range("AG2").Value = -24.63600000000001
$ws.Range("AH2").Value = 0.1966664535243089
$ws.Range("AI2").Value = 0.6578966131907308
$ws.Range("AJ2").Value = -0.0699361842254673
$ws.Range("AK2").Value = -1.055346127484579
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# Row 3
$ws.Range("D3").Value = -0.017
$ws.Range("K3").Value = -0.056
$ws.Range("L3").Value = -0.8484848484848485
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 0.256
$ws.Range("V3").Value = 0.01036437246963563
$ws.Range("W3").Value = -0.05333333333333333
$ws.Range("X3").Value = 0.03667539916908493
$ws.Range("Y3").Value = -0.09000873250241825
$ws.Range("Z3").Value = 0.08638743455497383
$ws.Range("AB3").Value = 0.03667539916908493
$ws.Range("AC3").Value = -0.03667539916908493
$ws.Range("AD3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -0.256
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.01047291768941254
$ws.Range("AK3").Value = -0.2091503267973856

# Row 4
$ws.Range("D4").Value = 0.0179
$ws.Range("E4").Value = 0.0258
$ws.Range("G4").Value = 0.005568181818181818
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6.08
$ws.Range("L4").Value = 0.3454545454545454
$ws.Range("M4").Value = 1.77
$ws.Range("N4").Value = 0.005253784505788068
$ws.Range("O4").Value = 0.2911184210526316
$ws.Range("P4").Value = 1.77
$ws.Range("Q4").Value = 0.005253784505788068
$ws.Range("R4").Value = 0.2911184210526316
$ws.Range("U4").Value = 112.9
$ws.Range("V4").Value = 0.3351142772336005
$ws.Range("W4").Value = 0.1804154302670623
$ws.Range("X4").Value = 0.03673057250252455
$ws.Range("Y4").Value = 0.1436848577645378
$ws.Range("Z4").Value = -0.124031007751938
$ws.Range("AA4").Value = -0
$ws.Range("AB4").Value = 0.03669751335046647
$ws.Range("AC4").Value = -0.03669751335046647
$ws.Range("AD4").Value = 0.97
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.97
$ws.Range("AG4").Value = -111.93
$ws.Range("AH4").Value = 0.002870926687779323
$ws.Range("AI4").Value = 0.02696691687517376
$ws.Range("AJ4").Value = -0.4975330044005868
$ws.Range("AK4").Value = 1.454959053685168
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# Row 5
$ws.Range("D5").Value = -0.365
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -3.58
$ws.Range("L5").Value = -4.00895856662934
$ws.Range("U5").Value = 3.75
$ws.Range("V5").Value = 0.2450980392156863
$ws.Range("W5").Value = -0.3630831643002029
$ws.Range("X5").Value = 0.1510258410577266
$ws.Range("Y5").Value = -0.5141090053579295
$ws.Range("Z5").Value = 0.009758602978941962
$ws.Range("AB5").Value = 0.04722954430850716
$ws.Range("AC5").Value = -0.04722954430850716
$ws.Range("AD5").Value = 91.3
$ws.Range("AF5").Value = 91.3
$ws.Range("AG5").Value = 87.55
$ws.Range("AH5").Value = 0.8564727954971858
$ws.Range("AI5").Value = 0.8881322957198443
$ws.Range("AJ5").Value = 0.8512396694214877
$ws.Range("AK5").Value = 0.883897021706209
